$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$k2Text = @'
Total number of records matching between DB & Response: 2, below are the test steps for this test case
'@
$i3Text = @'
Response_dependentRelationshipId: 5465997346542006717
DB_dependentRelationshipId: 5465997346542006717
Response_dependentRelationshipDescription: Test_Countryb2_To_Test_Statee3
DB_dependentRelationshipDescription: Test_Countryb2_To_Test_Statee3

'@
$a5Text = @'
TC_02
'@
$b5Text = @'
Verify the error message when passing the multiple(2 attributes) invalid attribute names
'@
$c5Text = @'
GraphQL
'@
$d5Text = @'

{
	"query":"
	{
		  dependentCountryRelationshipTypes 
		{
			    dependentRelationshipIds    dependentRelationshipDescriptions  
		}
	}"
}
'@
$f5Text = @'
NA
'@
$i5Text = @'

{
	"meta":
	{
		"version":"1.0.0",
		"errors":
		[
			{
				"error":"ValidationError",
				"message":"Validation error of type FieldUndefined: Field 'dependentRelationshipIds' in type 'DependentCountryRelationship' is undefined @ 'dependentCountryRelationshipTypes/dependentRelationshipIds'",
				"path":null
			},
			{
				"error":"ValidationError",
				"message":"Validation error of type FieldUndefined: Field 'dependentRelationshipDescriptions' in type 'DependentCountryRelationship' is undefined @ 'dependentCountryRelationshipTypes/dependentRelationshipDescriptions'",
				"path":null
			}
		]
	},
	"data":null
}
'@

# Update TC_01 summary comment (K2) and its I3 per-record detail to new counts/data
$ws.Range("K2").Value = $k2Text
$ws.Range("I3").Value = $i3Text

# Repurpose row 5 into the new TC_02 test case row
$ws.Range("A5").Value = $a5Text
$ws.Range("B5").Value = $b5Text
$ws.Range("C5").Value = $c5Text
$ws.Range("D5").Value = $d5Text
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = $f5Text
$ws.Range("H5").Value = "200"
$ws.Range("I5").Value = $i5Text

# Drop the two now-obsolete detail rows for TC_01 (rows 6 and 7)
$ws.Rows("6:7").Delete()

